# Update benchmark: 2025-11-07 06:38:56 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BENCHMARK")

$ws.Range("H2").Value = "15 TL - 15 TL"
$ws.Range("G3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("G4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("G5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("G6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("G8").Value = "15,23 TL - 30,47 TL - 211,05 TL"
$ws.Range("G9").Value = "15,23 TL - 30,47 TL - 211,05 TL"
$ws.Range("G10").Value = "15,23 TL - 30,47 TL - 211,05 TL"
$ws.Range("G11").Value = "3,04 TL - 6,09 TL - 76,17 TL"
$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 9.999.999.999.999 TL"
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 53,19 TL"
$ws.Range("G14").Value = "6.300 TL - 6,09 TL"
